$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 -- copy formatting (bold, centered, bordered) from H1, then set text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data cells I2:J56 (parallel flat arrays, one entry per row 2..56)
$iVals = @(6,3,4,5,6,7,8,3,6,7,10,9,7,6,6,3,6,3,8,8,10,7,8,8,6,8,5,4,6,6,6,6,6,7,8,6,6,6,9,9,7,3,8,5,9,2,5,10,9,1,8,1,8,5,6)
$jVals = @(6,3,4,5,6,7,8,4,6,7,10,9,7,6,6,3,6,3,8,9,10,7,8,8,6,9,6,5,6,6,6,6,7,7,9,6,6,6,9,9,8,4,8,5,9,3,6,11,9,1,8,1,8,6,6)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
